$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated plot data values for the balance service providing scenario
$data = @(
    @(1, 2104, 2515, 2611),
    @(2, 886.08300000000008, 1925.7290000000005, 2129.4700000000003),
    @(3, 11543.705000000004, 6488.0490000000018, 6737.002000000004),
    @(4, 6659.7100000000019, 3682.746000000001, 4209.5439999999999),
    @(5, 4494.0480000000007, 2247.5009999999997, 2209.3029999999999),
    @(6, 0, 0, 0),
    @(7, 3489, 2788, 2807),
    @(8, 769.35000000000014, 2061.8949999999991, 2227.2549999999992),
    @(9, 10928.412999999993, 6145.0419999999995, 6425.9440000000004),
    @(10, 6190.8819999999969, 3433.9150000000004, 4066.0949999999993),
    @(11, 4529.7630000000008, 2226.753999999999, 2279.1679999999992),
    @(12, 0, 0, 0),
    @(13, 3586, 2166, 1995),
    @(14, 1680.8300000000006, 2365.848, 2592.3199999999993),
    @(15, 10898.198999999999, 6496.1960000000036, 6838.2979999999998),
    @(16, 6418.5089999999991, 3538.2750000000015, 4137.1259999999975),
    @(17, 4330.2969999999987, 2277.2980000000002, 2393.0169999999994),
    @(18, 0, 0, 0),
    @(19, 1373.362499999999, 974.58799999999962, 904.50499999999965),
    @(20, 2226.7299999999991, 1618.7199999999998, 1595.3999999999992),
    @(21, 11257.601999999997, 5283.1869999999981, 5273.3289999999988),
    @(22, 6446.7100000000028, 2889.52, 3212.8600000000006),
    @(23, 4436.0550000000012, 1788.8079999999998, 1895.6379999999997),
    @(24, 0, 0, 0),
    @(25, 1437.6405000000002, 1305.6855000000005, 1296.9269999999997),
    @(26, 2168.4419999999991, 1919.04, 1950.3590000000004),
    @(27, 9275.2019999999975, 5731.1030000000019, 5661.070999999999),
    @(28, 5620.2959999999975, 3200.3049999999994, 3448.8949999999995),
    @(29, 3701.6649999999995, 1943.5560000000007, 1983.7219999999991),
    @(30, 0, 0, 0),
    @(31, 1512.3250000000003, 1711.5965000000003, 1658.5349999999992),
    @(32, 1903.3499999999999, 1415.1059999999993, 1425.0379999999998),
    @(33, 16518.094000000001, 10803.777, 10427.597999999998),
    @(34, 5041.8729999999978, 3122.7190000000005, 3308.4050000000007),
    @(35, 4795.8240000000023, 2852.4750000000004, 2852.4190000000003),
    @(36, 5.8399999999996908, 133.69000000000005, 121.76000000000067),
    @(37, 5059.6164999999946, 3905.7699999999977, 3782.1229999999991),
    @(38, 1969.172, 892.43800000000022, 945.06500000000028),
    @(39, 26598.50299999999, 14668.986999999999, 14325.897999999994),
    @(40, 5868.4619999999932, 3155.5130000000022, 3470.2359999999999),
    @(41, 8294.1510000000035, 4466.217999999998, 4292.9000000000015),
    @(42, 0, 0, 0),
    @(43, 5346.5099999999984, 4057.3584999999989, 4199.5489999999991),
    @(44, 1541.1400000000006, 774.31000000000063, 781.94999999999982),
    @(45, 23306.753000000026, 14927.280000000004, 14929.637999999992),
    @(46, 5499.7400000000034, 3429.6760000000004, 3680.5840000000017),
    @(47, 7589.8439999999946, 4775.7479999999987, 4709.6369999999997),
    @(48, 0, 0, 0),
    @(49, 3887.4975000000031, 2863.0859999999971, 3005.5939999999991),
    @(50, 2235.9100000000017, 1762.1320000000005, 1902.9730000000004),
    @(51, 16759.309999999994, 9248.1730000000025, 9299.146999999999),
    @(52, 4869.9800000000023, 2704.0450000000005, 2941.7550000000001),
    @(53, 5089.1450000000059, 2495.0159999999996, 2607.3530000000005),
    @(54, 0, 0, 0),
    @(55, 3381.9580000000005, 2380.5339999999983, 2428.5495000000019),
    @(56, 2284.5360000000001, 1891.0830000000001, 2068.4569999999999),
    @(57, 12252.739999999998, 6782.1049999999932, 6758.3839999999964),
    @(58, 4328.6750000000002, 2233.5149999999994, 2605.3999999999987),
    @(59, 3628.4829999999993, 1841.2869999999996, 1956.3139999999999),
    @(60, 0, 0, 0),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Column width changes: columns A and B narrow from 11.71 to 9.71
$ws.Columns.Item(1).ColumnWidth = 9.7109375
$ws.Columns.Item(2).ColumnWidth = 9.7109375
